$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins / Losses / Ties in AC1:AE1, matching the style of
# the existing header row (bold, bordered, centered - same as AB1).
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

# Team record (same W/L/T for every player row) for rows 2-42.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 29).Value = 74   # AC -> Wins
    $ws.Cells.Item($r, 30).Value = 88   # AD -> Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE -> Ties
}
